$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.077.05'
$ws.Range('E2').Value = '  +3.82%  '
$ws.Range('D3').Value = '2.648.44'
$ws.Range('E3').Value = '  +6.00%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '326.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.83%  '
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.557'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.00%  '
$ws.Range('D15').Value = '3.063.40'
$ws.Range('E15').Value = '  +6.17%  '
$ws.Range('D16').Value = '2.662.91'
$ws.Range('E16').Value = '  +6.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.872'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.09%  '
$ws.Range('D18').Value = '49.977.71'
$ws.Range('E18').Value = '  +4.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.57%  '
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +2.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '276.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  +3.50%  '
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.99%  '
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0815'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.38%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.52%  '
$ws.Range('E38').Value = '  +6.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.85%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '123.72'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.113'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('E42').Value = '  +0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0316'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.68%  '
$ws.Range('D45').Value = '2.084.02'
$ws.Range('E45').Value = '  +4.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +15.55%  '
$ws.Range('E48').Value = '  +4.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.15'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.40'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.96%  '
